$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", copying the existing header style (G1) so
# it reuses the same style index instead of creating a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column values (H2:H5), unstyled like the other numeric data cells.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
